$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.106.74"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "1.912.27"
$ws.Range("E3").Value = "  -2.83%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -1.46%  "
$ws.Range("D5").Value = "'327.71"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("D7").Value = "'0.4682"
$ws.Range("E7").Value = "  -5.44%  "
$ws.Range("D8").Value = "'0.4011"
$ws.Range("E8").Value = "  -3.87%  "
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("D10").Value = "'0.08381"
$ws.Range("E10").Value = "  -9.01%  "
$ws.Range("D11").Value = "'1.042"
$ws.Range("E11").Value = "  -4.34%  "
$ws.Range("D12").Value = "'22.09"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "1.974.48"
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").Value = "'7.438"
$ws.Range("E14").Value = "  -5.17%  "
$ws.Range("D15").Value = "'6.056"
$ws.Range("E15").Value = "  -5.74%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "'89.46"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "'0.00001059"
$ws.Range("E18").Value = "  -4.23%  "
$ws.Range("D19").Value = "'0.06580"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("D20").Value = "'17.94"
$ws.Range("E20").Value = "  -5.80%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "'5.714"
$ws.Range("E22").Value = "  -3.94%  "
$ws.Range("D23").Value = "28.128.67"
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("D24").Value = "'11.32"
$ws.Range("E24").Value = "  -4.85%  "
$ws.Range("D25").Value = "'2.273"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "2.190.28"
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("D27").Value = "'153.77"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").Value = "'19.97"
$ws.Range("E28").Value = "  -3.40%  "
$ws.Range("D29").Value = "'2.122"
$ws.Range("E29").Value = "  -5.63%  "
$ws.Range("D30").Value = "'5.660"
$ws.Range("E30").Value = "  -8.73%  "
$ws.Range("D31").Value = "'122.82"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").Value = "'0.9722"
$ws.Range("E32").Value = "  -6.41%  "
$ws.Range("D33").Value = "'0.09561"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("D34").Value = "'1.439"
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("D35").Value = "'3.636"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").Value = "'5.526"
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("D37").Value = "'8.794"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").Value = "'0.02294"
$ws.Range("E38").Value = "  -4.71%  "
$ws.Range("D39").Value = "'0.06137"
$ws.Range("E39").Value = "  -3.55%  "
$ws.Range("D40").Value = "'1.219"
$ws.Range("E40").Value = "  -7.28%  "
$ws.Range("D41").Value = "'0.6111"
$ws.Range("E41").Value = "  -5.06%  "
$ws.Range("D42").Value = "'10.99"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").Value = "'0.1896"
$ws.Range("E44").Value = "  -4.74%  "
$ws.Range("E45").Value = "  -3.26%  "
$ws.Range("D46").Value = "'0.5831"
$ws.Range("E46").Value = "  -5.24%  "
$ws.Range("D47").Value = "'12.64"
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("D48").Value = "'2.017"
$ws.Range("E48").Value = "  -6.50%  "
$ws.Range("D49").Value = "'3.448"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Value = "'0.06839"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").Value = "'1.080"
$ws.Range("E51").Value = "  -3.11%  "
